$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (17-dec) before the
#     existing "01-oct." column (ET), shifting everything after it one
#     column to the right (ET:FX -> EU:FY).
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns("ET").Insert()
$ws1.Range("ET1").Value = "17-dec"
$ws1.Range("ET2:ET25").Value = "-"

# --- Sheet "Gaz": append the new daily row.
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A180").NumberFormat = "@"
$ws2.Range("A180").Value = "2025-12-15"
$ws2.Range("A180").Style = "Normal"
$ws2.Range("B180").Value = 25.68

# --- Sheet "CO2": append the new daily row.
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A180").NumberFormat = "@"
$ws3.Range("A180").Value = "2025-12-15"
$ws3.Range("A180").Style = "Normal"
$ws3.Range("B180").Value = 84.59999999999999
